# Update version string from "mines - version 1.0.0 (Feb 3 2026) (built on
# February 03 2026 10.14.00 EST)" to
# "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on
# February 03 2026 17.29.55 EST)" throughout the workbook.

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Dendrobium Coal Mine, Australia, M0034, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($r = 2; $r -le 8; $r++) {
    $wsData.Range("S" + $r).Value = $newVersion
}
